$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pair (C3 -> Lrp1) specificity table: now covers the full 3x3
# sending/target cluster cross-product (ECs, FAPs, sCs) per Dr Hou advice,
# replacing the previous partial (2 target-cluster) table.
$rows = @(
    @("ECs", "C3", "Lrp1", "ECs", 2, 0.6666666666666666, 34.739995, 104.219985, 0.1827267341390226, 0.1827267341390226, 3, 1, 17.16653, 51.49959, 0.0560345397128279, 0.0560345397128279, 596.36516636735, 5367.28649730615, 0.01023900844070841, 0.01023900844070841),
    @("ECs", "C3", "Lrp1", "FAPs", 2, 0.6666666666666666, 34.739995, 104.219985, 0.1827267341390226, 0.1827267341390226, 3, 1, 256.4443053333333, 769.332916, 0.8370788162388805, 0.8370788162388805, 8908.873885058472, 80179.86496552626, 0.1529566783082897, 0.1529566783082897),
    @("ECs", "C3", "Lrp1", "sCs", 2, 0.6666666666666666, 34.739995, 104.219985, 0.1827267341390226, 0.1827267341390226, 3, 1, 32.74538866666666, 98.236166, 0.1068866440482915, 0.1068866440482915, 1137.574638553057, 10238.17174697751, 0.01953104739002451, 0.01953104739002451),
    @("FAPs", "C3", "Lrp1", "ECs", 3, 1, 154.8642143333334, 464.5926430000001, 0.8145606273154508, 0.8145606273154508, 3, 1, 17.16653, 51.49959, 0.0560345397128279, 0.0560345397128279, 2658.481181279597, 23926.33063151637, 0.04564352981981363, 0.04564352981981363),
    @("FAPs", "C3", "Lrp1", "FAPs", 3, 1, 154.8642143333334, 464.5926430000001, 0.8145606273154508, 0.8145606273154508, 3, 1, 256.4443053333333, 769.332916, 0.8370788162388805, 0.8370788162388805, 39714.04586570412, 357426.412791337, 0.6818514456680174, 0.6818514456680174),
    @("FAPs", "C3", "Lrp1", "sCs", 3, 1, 154.8642143333334, 464.5926430000001, 0.8145606273154508, 0.8145606273154508, 3, 1, 32.74538866666666, 98.236166, 0.1068866440482915, 0.1068866440482915, 5071.088888902971, 45639.80000012674, 0.08706565182761963, 0.08706565182761963),
    @("sCs", "C3", "Lrp1", "ECs", 3, 1, 0.5157266666666667, 1.54718, 0.002712638545526686, 0.002712638545526686, 3, 1, 17.16653, 51.49959, 0.0560345397128279, 0.0560345397128279, 8.853237295133333, 79.6791356562, 0.0001520014523058628, 0.0001520014523058628),
    @("sCs", "C3", "Lrp1", "FAPs", 3, 1, 0.5157266666666667, 1.54718, 0.002712638545526686, 0.002712638545526686, 3, 1, 256.4443053333333, 769.332916, 0.8370788162388805, 0.8370788162388805, 132.2551667752089, 1190.29650097688, 0.002270692262573436, 0.002270692262573437),
    @("sCs", "C3", "Lrp1", "sCs", 3, 1, 0.5157266666666667, 1.54718, 0.002712638545526686, 0.002712638545526686, 3, 1, 32.74538866666666, 98.236166, 0.1068866440482915, 0.1068866440482915, 16.88767014576444, 151.98903131188, 0.0002899448306473861, 0.0002899448306473861)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $rows[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
